# Update "want to go" counts (column F) on the exhibition, performance and
# "all types" sheets to the values generated at 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 603
$wsExpo.Range("F7").Value = 15041
$wsExpo.Range("F9").Value = 5
$wsExpo.Range("F11").Value = 15246
$wsExpo.Range("F12").Value = 38
$wsExpo.Range("F13").Value = 8770
$wsExpo.Range("F14").Value = 339
$wsExpo.Range("F16").Value = 69
$wsExpo.Range("F17").Value = 182
$wsExpo.Range("F19").Value = 184
$wsExpo.Range("F20").Value = 14
$wsExpo.Range("F21").Value = 25
$wsExpo.Range("F22").Value = 519
$wsExpo.Range("F26").Value = 1088
$wsExpo.Range("F29").Value = 57
$wsExpo.Range("F36").Value = 276
$wsExpo.Range("F37").Value = 430
$wsExpo.Range("F39").Value = 5388

# Sheet "演出" (Performances)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 1005

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 603
$wsAll.Range("F7").Value = 15041
$wsAll.Range("F9").Value = 5
$wsAll.Range("F11").Value = 15246
$wsAll.Range("F12").Value = 38
$wsAll.Range("F13").Value = 8770
$wsAll.Range("F14").Value = 339
$wsAll.Range("F16").Value = 1005
$wsAll.Range("F17").Value = 69
$wsAll.Range("F18").Value = 182
$wsAll.Range("F20").Value = 184
$wsAll.Range("F21").Value = 14
$wsAll.Range("F22").Value = 25
$wsAll.Range("F23").Value = 519
$wsAll.Range("F27").Value = 1088
$wsAll.Range("F30").Value = 57
$wsAll.Range("F39").Value = 276
$wsAll.Range("F40").Value = 430
$wsAll.Range("F42").Value = 5388
